# Auto-generated edit script: updates Leve profit-calculation columns (H:N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets to reflect refreshed
# Universalis market-price data (scheduled runner sync).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 11 (Leve Item ID 5533)
$ws.Range("H11").Value = 1204.8182
$ws.Range("I11").Value = 1204.8182
$ws.Range("K11").Value = 1204.8182
$ws.Range("M11").Value = -1064.8182

# row 55 (Leve Item ID 5517)
$ws.Range("H55").Value = 58745.26
$ws.Range("I55").Value = 319.2
$ws.Range("J55").Value = 79611.71000000001
$ws.Range("K55").Value = 319.2
$ws.Range("L55").Value = 79611.71000000001
$ws.Range("M55").Value = -105.2
$ws.Range("N55").Value = -80039.71000000001

# row 74 (Leve Item ID 5507)
$ws.Range("H74").Value = 7284.3335
$ws.Range("I74").Value = 6419.5264
$ws.Range("K74").Value = 6419.5264
$ws.Range("M74").Value = -5483.5264

# row 77 (Leve Item ID 5507)
$ws.Range("H77").Value = 7284.3335
$ws.Range("I77").Value = 6419.5264
$ws.Range("K77").Value = 32097.632
$ws.Range("M77").Value = -27417.632

# row 87 (Leve Item ID 10651)
$ws.Range("H87").Value = 74807.5
$ws.Range("J87").Value = 74807.5
$ws.Range("L87").Value = 74807.5
$ws.Range("N87").Value = -77303.5

# row 90 (Leve Item ID 10651)
$ws.Range("H90").Value = 74807.5
$ws.Range("J90").Value = 74807.5
$ws.Range("L90").Value = 224422.5
$ws.Range("N90").Value = -236902.5

# row 99 (Leve Item ID 19883)
$ws.Range("H99").Value = 126559.625
$ws.Range("J99").Value = 168553.17
$ws.Range("L99").Value = 505659.51
$ws.Range("N99").Value = -508655.51

# row 113 (Leve Item ID 27775)
$ws.Range("H113").Value = 7371.636
$ws.Range("I113").Value = 7371.636
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 7371.636
$ws.Range("L113").Value = 0
$ws.Range("N113").Value = -4117.636
$ws.Range("M113").ClearContents()

# row 132 (Leve Item ID 44049)
$ws.Range("H132").Value = 3455.4736
$ws.Range("I132").Value = 3395.5095
$ws.Range("J132").Value = 4250
$ws.Range("K132").Value = 10186.5285
$ws.Range("L132").Value = 12750
$ws.Range("M132").Value = -7656.5285
$ws.Range("N132").Value = -17810

# row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 7848.9546
$ws.Range("I137").Value = 1317.6786
$ws.Range("J137").Value = 12661.474
$ws.Range("K137").Value = 3953.0358
$ws.Range("L137").Value = 37984.422
$ws.Range("M137").Value = -1403.0358
$ws.Range("N137").Value = -43084.422

$ws = $wb.Worksheets.Item("ARM")
# row 45 (Leve Item ID 27714)
$ws.Range("H45").Value = 93873.63
$ws.Range("I45").Value = 169151.67
$ws.Range("K45").Value = 169151.67
$ws.Range("M45").Value = -168774.67

# row 92 (Leve Item ID 18050)
$ws.Range("H92").Value = 30031938
$ws.Range("J92").Value = 25039922
$ws.Range("L92").Value = 25039922
$ws.Range("N92").Value = -25044914

# row 102 (Leve Item ID 19945)
$ws.Range("H102").Value = 1827
$ws.Range("I102").Value = 1532.6875
$ws.Range("K102").Value = 1532.6875
$ws.Range("M102").Value = 89.3125

# row 122 (Leve Item ID 36168)
$ws.Range("H122").Value = 41670492
$ws.Range("I122").Value = 76924570
$ws.Range("J122").Value = 6586.091
$ws.Range("K122").Value = 230773710
$ws.Range("L122").Value = 19758.273
$ws.Range("M122").Value = -230771260
$ws.Range("N122").Value = -24658.273

$ws = $wb.Worksheets.Item("BSM")
# row 86 (Leve Item ID 12526)
$ws.Range("H86").Value = 2431657.8
$ws.Range("I86").Value = 3402579.5
$ws.Range("J86").Value = 4353.5
$ws.Range("K86").Value = 3402579.5
$ws.Range("L86").Value = 4353.5
$ws.Range("M86").Value = -3401456.5
$ws.Range("N86").Value = -6599.5

# row 88 (Leve Item ID 10626)
$ws.Range("H88").Value = 30343
$ws.Range("J88").Value = 30343
$ws.Range("L88").Value = 30343
$ws.Range("N88").Value = -31155

# row 89 (Leve Item ID 12526)
$ws.Range("H89").Value = 2431657.8
$ws.Range("I89").Value = 3402579.5
$ws.Range("J89").Value = 4353.5
$ws.Range("K89").Value = 17012897.5
$ws.Range("L89").Value = 21767.5
$ws.Range("M89").Value = -17007281.5
$ws.Range("N89").Value = -32999.5

# row 91 (Leve Item ID 10626)
$ws.Range("H91").Value = 30343
$ws.Range("J91").Value = 30343
$ws.Range("L91").Value = 30343
$ws.Range("N91").Value = -33151

# row 105 (Leve Item ID 19947)
$ws.Range("H105").Value = 15154387
$ws.Range("I105").Value = 1259.25
$ws.Range("J105").Value = 23813316
$ws.Range("K105").Value = 1259.25
$ws.Range("L105").Value = 23813316
$ws.Range("M105").Value = 487.75
$ws.Range("N105").Value = -23816810

$ws = $wb.Worksheets.Item("CRP")
# row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 497998.8
$ws.Range("I31").Value = 8475471
$ws.Range("J31").Value = 14515.637
$ws.Range("K31").Value = 8475471
$ws.Range("L31").Value = 14515.637
$ws.Range("M31").Value = -8475176
$ws.Range("N31").Value = -15105.637

# row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 497998.8
$ws.Range("I34").Value = 8475471
$ws.Range("J34").Value = 14515.637
$ws.Range("K34").Value = 8475471
$ws.Range("L34").Value = 14515.637
$ws.Range("M34").Value = -8475269
$ws.Range("N34").Value = -14919.637

# row 99 (Leve Item ID 36198)
$ws.Range("H99").Value = 6086.4546
$ws.Range("I99").Value = 3488
$ws.Range("K99").Value = 3488
$ws.Range("M99").Value = -1990

# row 105 (Leve Item ID 19928)
$ws.Range("H105").Value = 2793.182
$ws.Range("I105").Value = 2314.4443
$ws.Range("K105").Value = 2314.4443
$ws.Range("M105").Value = -567.4443000000001

# row 126 (Leve Item ID 36198)
$ws.Range("H126").Value = 6086.4546
$ws.Range("I126").Value = 3488
$ws.Range("K126").Value = 10464
$ws.Range("M126").Value = -7994

$ws = $wb.Worksheets.Item("CUL")
# row 44 (Leve Item ID 4702)
$ws.Range("H44").Value = 999
$ws.Range("I44").Value = 999
$ws.Range("K44").Value = 2997
$ws.Range("M44").Value = -2599

# row 74 (Leve Item ID 12859)
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()

# row 77 (Leve Item ID 12859)
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# row 70 (Leve Item ID 14146)
$ws.Range("H70").Value = 6509.647
$ws.Range("I70").Value = 6106.273
$ws.Range("K70").Value = 6106.273
$ws.Range("M70").Value = -5836.273

# row 73 (Leve Item ID 14146)
$ws.Range("H73").Value = 6509.647
$ws.Range("I73").Value = 6106.273
$ws.Range("K73").Value = 6106.273
$ws.Range("M73").Value = -5170.273

# row 102 (Leve Item ID 36169)
$ws.Range("H102").Value = 14446.827
$ws.Range("I102").Value = 18045.857
$ws.Range("K102").Value = 18045.857
$ws.Range("M102").Value = -16423.857

# row 122 (Leve Item ID 36182)
$ws.Range("H122").Value = 413101.44
$ws.Range("I122").Value = 505124.5
$ws.Range("K122").Value = 1515373.5
$ws.Range("M122").Value = -1512923.5

# row 135 (Leve Item ID 42006)
$ws.Range("H135").Value = 112500
$ws.Range("J135").Value = 112500
$ws.Range("L135").Value = 112500
$ws.Range("N135").Value = -122640

$ws = $wb.Worksheets.Item("LTW")
# row 16 (Leve Item ID 5289)
$ws.Range("H16").Value = 1720.9584
$ws.Range("I16").Value = 1750.2273
$ws.Range("K16").Value = 1750.2273
$ws.Range("M16").Value = -1580.2273

# row 22 (Leve Item ID 5277)
$ws.Range("H22").Value = 2250
$ws.Range("I22").Value = 2512.75
$ws.Range("J22").Value = 1987.25
$ws.Range("K22").Value = 2512.75
$ws.Range("L22").Value = 1987.25
$ws.Range("M22").Value = -2217.75
$ws.Range("N22").Value = -2577.25

# row 27 (Leve Item ID 5277)
$ws.Range("H27").Value = 2250
$ws.Range("I27").Value = 2512.75
$ws.Range("J27").Value = 1987.25
$ws.Range("K27").Value = 2512.75
$ws.Range("L27").Value = 1987.25
$ws.Range("M27").Value = -2405.75
$ws.Range("N27").Value = -2201.25

# row 40 (Leve Item ID 36248)
$ws.Range("H40").Value = 914869.6
$ws.Range("I40").Value = 1254457.6
$ws.Range("K40").Value = 1254457.6
$ws.Range("M40").Value = -1254321.6

# row 46 (Leve Item ID 5282)
$ws.Range("H46").Value = 3306
$ws.Range("J46").Value = 3487.5625
$ws.Range("L46").Value = 3487.5625
$ws.Range("N46").Value = -3863.5625

# row 82 (Leve Item ID 12565)
$ws.Range("H82").Value = 1882.5333
$ws.Range("I82").Value = 1516.75
$ws.Range("K82").Value = 1516.75
$ws.Range("M82").Value = -1155.75

# row 85 (Leve Item ID 12565)
$ws.Range("H85").Value = 1882.5333
$ws.Range("I85").Value = 1516.75
$ws.Range("K85").Value = 1516.75
$ws.Range("M85").Value = -268.75

# row 93 (Leve Item ID 19993)
$ws.Range("H93").Value = 3798.111
$ws.Range("I93").Value = 3888.8333
$ws.Range("J93").Value = 3616.6667
$ws.Range("K93").Value = 3888.8333
$ws.Range("L93").Value = 3616.6667
$ws.Range("M93").Value = -2640.8333
$ws.Range("N93").Value = -6112.6667

# row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 3348.25
$ws.Range("I132").Value = 2562.2195
$ws.Range("J132").Value = 7952.143
$ws.Range("K132").Value = 7686.6585
$ws.Range("L132").Value = 23856.429
$ws.Range("M132").Value = -5156.6585
$ws.Range("N132").Value = -28916.429

$ws = $wb.Worksheets.Item("WVR")
# row 58 (Leve Item ID 3187)
$ws.Range("H58").Value = 10000
$ws.Range("I58").Value = 10000
$ws.Range("K58").Value = 10000
$ws.Range("M58").Value = -9692

# row 122 (Leve Item ID 36208)
$ws.Range("H122").Value = 4130.1797
$ws.Range("I122").Value = 3521.1875
$ws.Range("K122").Value = 10563.5625
$ws.Range("M122").Value = -8113.5625

# row 126 (Leve Item ID 36210)
$ws.Range("H126").Value = 10026.375
$ws.Range("J126").Value = 11876.917
$ws.Range("L126").Value = 35630.751
$ws.Range("N126").Value = -40570.751

# row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 34746.8
$ws.Range("I132").Value = 2157.5217
$ws.Range("K132").Value = 6472.5651
$ws.Range("M132").Value = -3942.5651
